$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D7"  = -7.529999999999999
    "A9"  = -21.778
    "D12" = -7.269
    "A18" = -22.387
    "A20" = -20.2
    "D26" = -8.123000000000001
    "A27" = -21.919
    "D27" = -8.599
    "D29" = -7.292
    "D37" = -7.822
    "D38" = -7.722
    "D51" = -8.638000000000002
    "D55" = -7.918000000000001
    "A69" = -21.627
    "D69" = -7.220999999999999
    "D70" = -7.145000000000001
    "A76" = -20.043
    "A82" = -22.152
    "D83" = -8.373000000000001
    "D102" = -7.833
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
